# Updated cryptos list on Tue Sep 19 11:08:57 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume refresh to the cryptos worksheet, and
# re-orders a few rows where new coins (MXToken, Cronos, EnergySwap,
# Algorand, Mantle) bumped existing ones down the rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes a string into a cell, forcing text storage so numeric-looking
    # strings (e.g. "218.48") are not silently reinterpreted as numbers -
    # matches the source data, which stores every Price/Volume cell as text.
    param(
        [string]$CellRef,
        [string]$Value
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.230.49"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.648.31"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "218.48"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6 - XRP
Set-TextValue "D6" "0.510"
$ws.Range("E6").Value = "  +1.68%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.257"
$ws.Range("E8").Value = "  +0.86%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.06%  "

# Row 10 - Solana
Set-TextValue "D10" "20.24"
$ws.Range("E10").Value = "  +2.79%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.03%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.879.60"
$ws.Range("E12").Value = "  -0.15%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.654.77"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.53%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.82%  "

# Row 16 - Litecoin
Set-TextValue "D16" "67.92"
$ws.Range("E16").Value = "  +2.50%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.209.40"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.20%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "220.92"
$ws.Range("E19").Value = "  -0.75%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.10%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -1.10%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +3.13%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -0.43%  "

# Row 25 - Monero
Set-TextValue "D25" "148.20"

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.07%  "

# Row 27 - Cosmos
Set-TextValue "D27" "7.41"
$ws.Range("E27").Value = "  +0.53%  "

# Row 28 - Stellar
Set-TextValue "D28" "0.119"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.82"
$ws.Range("E29").Value = "  -0.85%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.64%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.62%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.87%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.04"
$ws.Range("E33").Value = "  +0.35%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -0.10%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.274.45"
$ws.Range("E35").Value = "  +0.62%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +1.01%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.543"
$ws.Range("E38").Value = "  +0.44%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +1.95%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.07%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  +0.24%  "

# Row 42 - was FraxShare, now MXToken
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.22"
$ws.Range("E42").Value = "  +7.82%  "

# Row 43 - was MXToken, now FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.41"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.790.46"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - Aave
Set-TextValue "D45" "63.12"
$ws.Range("E45").Value = "  +1.65%  "

# Row 46 - Quant
Set-TextValue "D46" "92.60"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -1.36%  "

# Row 48 - was BabyDogeCoin, now Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0514"
$ws.Range("E48").Value = "  -0.51%  "

# Row 49 - was Cronos, now EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "7.72"
$ws.Range("E49").Value = "  +0.97%  "

# Row 50 - was EnergySwap, now Algorand
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0976"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51 - was Algorand, now Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.406"
$ws.Range("E51").Value = "  +0.08%  "
